# Updates the cryptos price/volume table with freshly scraped values
# (commit: "Updated symbol list on Fri Jan 27 07:35:50 UTC 2023 with GitHub Actions").
# Column B/C: coin name / coinranking.com link (plain text).
# Column D/E: price / 1h volume change, stored as TEXT in the sheet even though
# the strings look numeric (e.g. "305.23", "-0.77%") - so numeric-looking values
# are entered with a leading apostrophe to force text entry, matching the
# source workbook's inlineStr cells, then the cell style is reset to "Normal"
# since Excel tags quote-prefixed entries with a distinct cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = [ordered]@{
    "D2" = '305.23'
    "E2" = '-0.77%'
    "D3" = '35.76'
    "E3" = '-0.38%'
    "D4" = '5.039'
    "E4" = '-1.38%'
    "D5" = '0.07962'
    "E5" = '-1.59%'
    "D6" = '1.909'
    "E6" = '-2.26%'
    "B7" = 'GateToken'
    "C7" = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    "D7" = '4.138'
    "E7" = '-1.45%'
    "B8" = 'KuCoinToken'
    "C8" = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
    "D8" = '7.773'
    "E8" = '0.20%'
    "B9" = 'MXToken'
    "C9" = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    "D9" = '0.9194'
    "E9" = '-1.01%'
    "B10" = 'LiechtensteinCryptoassetsExchange'
    "C10" = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    "D10" = '0.1309'
    "E10" = '-4.69%'
    "B11" = 'WazirX'
    "C11" = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    "D11" = '0.1911'
    "E11" = '-0.41%'
    "B12" = 'MandalaExchangeToken'
    "C12" = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    "D12" = '0.09074'
    "E12" = '-1.57%'
    "B13" = 'BitrueCoin'
    "C13" = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    "D13" = '0.03439'
    "E13" = '-0.43%'
    "B14" = 'BitMartToken'
    "C14" = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    "D14" = '0.09838'
    "E14" = '-0.07%'
    "B15" = 'BitForexToken'
    "C15" = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    "D15" = '0.001403'
    "E15" = '-1.19%'
    "B16" = 'TigerCash'
    "C16" = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    "D16" = '0.006167'
    "E16" = '6.20%'
    "B17" = 'LEO'
    "C17" = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    "D17" = '3.714'
    "E17" = '2.62%'
    "E18" = '12.13%'
    "D19" = '0.3444'
    "E19" = '-0.07%'
    "D20" = '0.1341'
    "E20" = '2.03%'
    "D21" = '5.174'
    "E21" = '5.37%'
    "D22" = '0.2352'
    "E22" = '-3.81%'
    "D23" = '0.04412'
    "E23" = '-0.93%'
    "D24" = '0.001234'
    "E24" = '1.12%'
    "D25" = '0.004630'
    "E25" = '-4.14%'
    "D26" = '0.0001251'
    "E26" = '0.65%'
    "D27" = '0.0004445'
    "E27" = '0.07%'
    "D39" = '0.01938'
    "E39" = '-3.96%'
    "D40" = '0.05280'
    "E40" = '7.10%'
    "D41" = '0.007603'
    "E41" = '-1.00%'
    "D42" = '0.01013'
    "E42" = '-0.28%'
    "D43" = '0.1353'
    "E43" = '-1.67%'
    "D44" = '0.002152'
    "E44" = '2.23%'
    "D45" = '0.009949'
    "E45" = '-14.25%'
    "D46" = '0.00006114'
    "E46" = '-5.27%'
    "E47" = '-0.12%'
    "D48" = '63.57'
    "E48" = '0.01%'
    "D49" = '0.001659'
    "E49" = '39.11%'
    "D50" = '0.00002102'
    "E50" = '-0.12%'
    "D51" = '0.0002002'
    "E51" = '-0.12%'
}

foreach ($addr in $updates.Keys) {
    $newValue = $updates[$addr]
    $cell = $ws.Range($addr)

    $looksNumeric = $newValue -match '^-?[0-9]' 
    if ($looksNumeric) {
        $cell.Value = "'" + $newValue
        $cell.Style = "Normal"
    } else {
        $cell.Value = $newValue
    }
}
